# 🔄 MAJ automatique BRVM via GitHub Actions
#
# Sheet "Recommandations": refresh the sector/stock rows (2..34) with new
# figures and insert 5 new rows (so the sheet grows from 34 to 39 data
# rows), reordering several stock lines in the process.
#
# Sheet "Top_YTD": drop everything except the header row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: Recommandations
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Recommandations")

# Final state for rows 2..39 (Titre, Hausse, Baisse, VarTotale, DerniereVar, Recommandation, Strategie)
$rows = @(
    @('BRVM - SERVICES PUBLICS', 0, 8, 3276.37, 105.9, '🟡 Observer', '➖ Neutre'),
    @('BRVM - AUTRES SECTEURS', 0, 4, 2771.51, 669.2, '🟡 Observer', '➖ Neutre'),
    @('BRVM - DISTRIBUTION', 0, 4, 1985.75, 494.78, '🟡 Observer', '➖ Neutre'),
    @('BRVM - TRANSPORT', 0, 4, 1455.18, 370.83, '🟡 Observer', '➖ Neutre'),
    @('BRVM - AGRICULTURE', 0, 4, 1336.97, 341.63, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIE  (**)', 0, 4, 1024.6, 257.22, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRINCIPAL  (**)', 0, 4, 866.29, 217.65, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DE BASE  (**)', 0, 4, 856.07, 216.08, '🟡 Observer', '➖ Neutre'),
    @('BRVM - CONSOMMATION DISCRETIONNAIRE', 0, 4, 700.63, 173.02, '🟡 Observer', '➖ Neutre'),
    @('BRVM - FINANCES', 0, 4, 587.99, 148.41, '🟡 Observer', '➖ Neutre'),
    @('BRVM - SERVICES FINANCIERS', 0, 4, 577.87, 145.86, '🟡 Observer', '➖ Neutre'),
    @('BRVM-PRESTIGE', 0, 4, 568.16, 144.25, '🟡 Observer', '➖ Neutre'),
    @('BRVM - INDUSTRIELS', 0, 4, 529.82, 132.33, '🟡 Observer', '➖ Neutre'),
    @('BRVM - ENERGIE', 0, 4, 442.1, 110.91, '🟡 Observer', '➖ Neutre'),
    @('BRVM - TELECOMMUNICATIONS', 0, 4, 372.94, 94.9, '🟡 Observer', '➖ Neutre'),
    @('NEI-CEDA CI (NEIC)', 3, 0, 20.06, 7.21, '🟢 Achat', '✅ Renforcer'),
    @('SONATEL SN (SNTS)', 1, 0, 7.05, 7.05, '🟡 Observer', '➖ Neutre'),
    @('SAPH CI (SPHC)', 1, 0, 5.33, 5.33, '🟡 Observer', '➖ Neutre'),
    @('SOCIETE GENERALE COTE D''IVOIRE (SGBC)', 1, 0, 4.75, 4.75, '🟡 Observer', '➖ Neutre'),
    @('CIE CI (CIEC)', 1, 0, 4.66, 4.66, '🟡 Observer', '➖ Neutre'),
    @('BANK OF AFRICA NG (BOAN)', 1, 0, 4.13, 4.13, '🟡 Observer', '➖ Neutre'),
    @('VIVO ENERGY CI (SHEC)', 1, 0, 3.7, 3.7, '🟡 Observer', '➖ Neutre'),
    @('AFRICA GLOBAL LOGISTICS CI (SDSC)', 1, 0, 3.42, 3.42, '🟡 Observer', '➖ Neutre'),
    @('LOTERIE NATIONALE DU BENIN (LNBB)', 1, 1, 2.9, -3.65, '🟡 Observer', '👀 À surveiller'),
    @('SOCIETE IVOIRIENNE DE BANQUE  (SIBC)', 1, 1, 1.73, 3.23, '🟡 Observer', '👀 À surveiller'),
    @('SICOR CI (SICC)', 1, 1, 1, -5.71, '🟡 Observer', '👀 À surveiller'),
    @('TRACTAFRIC MOTORS CI (PRSC)', 1, 1, 0.81, -3.43, '🟡 Observer', '👀 À surveiller'),
    @('SETAO CI (STAC)', 2, 2, 0.63, -6.69, '🟡 Observer', '👀 À surveiller'),
    @('ECOBANK TRANS. INCORP. TG (ETIT)', 1, 1, 0.2, 4.55, '🟡 Observer', '👀 À surveiller'),
    @('TOTAL', 0, 3, 0, 0, '🟡 Observer', '➖ Neutre'),
    @('BANK OF AFRICA BF (BOABF)', 1, 1, -1.16, -6.02, '🟡 Observer', '👀 À surveiller'),
    @('CFAO MOTORS CI (CFAC)', 0, 1, -1.38, -1.38, '🟡 Observer', '➖ Neutre'),
    @('BANK OF AFRICA ML (BOAM)', 0, 1, -3.5, -3.5, '🟡 Observer', '➖ Neutre'),
    @('NESTLE CI (NTLC)', 0, 2, -3.74, -2.78, '🟡 Observer', '➖ Neutre'),
    @('ORAGROUP TOGO (ORGT)', 1, 2, -3.85, -3.81, '🟡 Observer', '👀 À surveiller'),
    @('SERVAIR ABIDJAN CI (ABJC)', 0, 1, -4.92, -4.92, '🟡 Observer', '➖ Neutre'),
    @('SUCRIVOIRE (SCRC)', 0, 2, -6.35, -5.45, '🟡 Observer', '➖ Neutre'),
    @('ERIUM CI (Ex AIR LIQUIDE CI) (SIVC)', 1, 3, -11, 5.33, '🔴 Vente', '⚠️ Risque de décrochage')
)

$r = 2
foreach ($row in $rows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet 2: Top_YTD - keep only the header row, drop rows 2..11
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top_YTD")
$ws2.Range("A2:B11").EntireRow.Delete()
